$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "60.323.60"
Set-TextValue $ws.Range("E2") "  +3.41%  "
Set-TextValue $ws.Range("D3") "3.219.74"
Set-TextValue $ws.Range("E3") "  +1.99%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "539.58"
Set-TextValue $ws.Range("E5") "  +0.39%  "
Set-TextValue $ws.Range("D6") "146.36"
Set-TextValue $ws.Range("E6") "  +4.54%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.06%  "
Set-TextValue $ws.Range("E8") "  +4.14%  "
Set-TextValue $ws.Range("D9") "7.38"
Set-TextValue $ws.Range("E9") "  +0.35%  "
Set-TextValue $ws.Range("E10") "  +4.17%  "
Set-TextValue $ws.Range("E11") "  +3.39%  "
Set-TextValue $ws.Range("D12") "3.766.72"
Set-TextValue $ws.Range("E12") "  +1.98%  "
Set-TextValue $ws.Range("D13") "0.139"
Set-TextValue $ws.Range("E13") "  -1.14%  "
Set-TextValue $ws.Range("D14") "26.36"
Set-TextValue $ws.Range("E14") "  +0.93%  "
Set-TextValue $ws.Range("D15") "0.0000175"
Set-TextValue $ws.Range("E15") "  +3.08%  "
Set-TextValue $ws.Range("D16") "60.331.20"
Set-TextValue $ws.Range("E16") "  +3.38%  "
Set-TextValue $ws.Range("D17") "3.196.43"
Set-TextValue $ws.Range("E17") "  +1.75%  "
Set-TextValue $ws.Range("D18") "6.28"
Set-TextValue $ws.Range("E18") "  +0.74%  "
Set-TextValue $ws.Range("D19") "13.24"
Set-TextValue $ws.Range("E19") "  +1.63%  "
Set-TextValue $ws.Range("D20") "8.39"
Set-TextValue $ws.Range("E20") "  +2.35%  "
Set-TextValue $ws.Range("D21") "381.86"
Set-TextValue $ws.Range("E21") "  +1.30%  "
Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.28%  "
Set-TextValue $ws.Range("D23") "0.530"
Set-TextValue $ws.Range("E23") "  +2.21%  "
Set-TextValue $ws.Range("D24") "70.37"
Set-TextValue $ws.Range("E24") "  -0.23%  "
Set-TextValue $ws.Range("D25") "8.95"
Set-TextValue $ws.Range("E25") "  +9.96%  "
Set-TextValue $ws.Range("D26") "0.171"
Set-TextValue $ws.Range("E26") "  +1.78%  "
Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  -0.25%  "
Set-TextValue $ws.Range("D28") "0.0₃0914"
Set-TextValue $ws.Range("E28") "  +3.37%  "
Set-TextValue $ws.Range("E29") "  +0.68%  "
Set-TextValue $ws.Range("B30") "RenderToken"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D30") "6.23"
Set-TextValue $ws.Range("E30") "  +0.66%  "
Set-TextValue $ws.Range("D31") "5.48"
Set-TextValue $ws.Range("E31") "  +5.66%  "
Set-TextValue $ws.Range("B32") "EthereumClassic"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D32") "22.49"
Set-TextValue $ws.Range("E32") "  +3.01%  "
Set-TextValue $ws.Range("E33") "  +3.97%  "
Set-TextValue $ws.Range("D34") "6.67"
Set-TextValue $ws.Range("E34") "  +6.55%  "
Set-TextValue $ws.Range("D35") "157.10"
Set-TextValue $ws.Range("E35") "  -2.67%  "
Set-TextValue $ws.Range("E36") "  +0.82%  "
Set-TextValue $ws.Range("D37") "2.799.67"
Set-TextValue $ws.Range("E37") "  +5.76%  "
Set-TextValue $ws.Range("D38") "25.96"
Set-TextValue $ws.Range("E38") "  +1.64%  "
Set-TextValue $ws.Range("D39") "0.0711"
Set-TextValue $ws.Range("E39") "  +4.31%  "
Set-TextValue $ws.Range("D40") "1.69"
Set-TextValue $ws.Range("E40") "  +0.22%  "
Set-TextValue $ws.Range("D41") "4.27"
Set-TextValue $ws.Range("E41") "  +0.50%  "
Set-TextValue $ws.Range("D42") "40.19"
Set-TextValue $ws.Range("E42") "  +4.08%  "
Set-TextValue $ws.Range("D43") "0.723"
Set-TextValue $ws.Range("E43") "  +2.66%  "
Set-TextValue $ws.Range("D44") "0.0288"
Set-TextValue $ws.Range("E44") "  +4.18%  "
Set-TextValue $ws.Range("D45") "3.254.60"
Set-TextValue $ws.Range("E45") "  +1.92%  "
Set-TextValue $ws.Range("B46") "ONDO"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D46") "1.01"
Set-TextValue $ws.Range("E46") "  +2.41%  "
Set-TextValue $ws.Range("B47") "Stellar"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D47") "0.104"
Set-TextValue $ws.Range("E47") "  +0.29%  "
Set-TextValue $ws.Range("D48") "6.19"
Set-TextValue $ws.Range("B49") "InjectiveProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D49") "20.91"
Set-TextValue $ws.Range("E49") "  +2.93%  "
Set-TextValue $ws.Range("B50") "SuiNetwork"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D50") "0.807"
Set-TextValue $ws.Range("E50") "  +6.94%  "
Set-TextValue $ws.Range("B51") "Bittensor"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D51") "273.20"
Set-TextValue $ws.Range("E51") "  +9.93%  "
